# "Turned on all checks." — the Workflow checklist had a single check
# ("Unused variables", row 17) left disabled ("No"), and its check file had
# moved into its own subfolder. Enable the check and point it at the new
# checker path, matching every other row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

$ws.Range("A17").Value = "Yes"
$ws.Range("C17").Value = "Checks\Standard\UnusedVariables\UnusedVariables.xaml"

# Column A (Run Yes/No) and column E (Action Fix/Double check) validations
# now cover one uninterrupted block (A2:A21 / E2:E21) instead of two
# contiguous pieces, since every row shares the same "Yes, No" / "Fix,
# Double check" list rule. Re-apply the same rule across the merged range.
$rngA = $ws.Range("A2:A21")
$rngA.Validation.Delete()
$rngA.Validation.Add(3, 1, 1, '"Yes, No"')
$rngA.Validation.IgnoreBlank = $true
$rngA.Validation.InCellDropdown = $true
$rngA.Validation.ShowInput = $true
$rngA.Validation.ShowError = $true

$rngE = $ws.Range("E2:E21")
$rngE.Validation.Delete()
$rngE.Validation.Add(3, 1, 1, '"Fix, Double check"')
$rngE.Validation.IgnoreBlank = $true
$rngE.Validation.InCellDropdown = $true
$rngE.Validation.ShowInput = $true
$rngE.Validation.ShowError = $true
